$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("35:36").Insert()

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 44914
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100103
$ws.Range("H35").Value = "Frutos de hueso (carozo)"
$ws.Range("I35").Value = 100103001
$ws.Range("J35").Value = "Cereza"
$ws.Range("K35").Value = "Santina"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 120
$ws.Range("N35").Value = 5000
$ws.Range("O35").Value = 5500
$ws.Range("P35").Value = 5250
$ws.Range("Q35").Value = "$/bandeja 10 kilos"
$ws.Range("R35").Value = "Provincia de Curicó"
$ws.Range("S35").Value = 525
$ws.Range("T35").Value = 10

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 44914
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100103
$ws.Range("H36").Value = "Frutos de hueso (carozo)"
$ws.Range("I36").Value = 100103001
$ws.Range("J36").Value = "Cereza"
$ws.Range("K36").Value = "Santina"
$ws.Range("L36").Value = "Segunda"
$ws.Range("M36").Value = 60
$ws.Range("N36").Value = 4000
$ws.Range("O36").Value = 4000
$ws.Range("P36").Value = 4000
$ws.Range("Q36").Value = "$/bandeja 10 kilos"
$ws.Range("R36").Value = "Provincia de Curicó"
$ws.Range("S36").Value = 400
$ws.Range("T36").Value = 10
